$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" for the 3d085b3c row (row 3)
$wsOverview.Range("G3").Value = "2016-08-22 22:44:54"

# zh-cn sheet: Correspond Handoff Datetime / Correspond Handback DateTime for 3d085b3c row (row 3)
$wsZhCn.Range("H3").Value = "2016-08-22 22:44:49"
$wsZhCn.Range("K3").Value = "2016-08-22 22:45:18"

# de-de sheet: Correspond Handback DateTime for 3d085b3c row (row 3)
$wsDeDe.Range("K3").Value = "2016-08-22 22:45:25"
